$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = New-Object "object[,]" 17,6

$rows[0,0] = 'cl:10147'
$rows[0,1] = 'data type'
$rows[0,2] = 'Data type is an attribute associated with a piece of data that tells a computer system how to interpret its value.'
$rows[0,3] = ''
$rows[0,4] = 'https://dataled.academy/guides/data-types/'
$rows[0,5] = ''

$rows[1,0] = 'cl:10148'
$rows[1,1] = 'integer'
$rows[1,2] = 'It is the most common numeric data type used to store numbers without a fractional component (-707, 0, 707).'
$rows[1,3] = 'int'
$rows[1,4] = 'https://dataled.academy/guides/data-types/'
$rows[1,5] = 'cl:10147'

$rows[2,0] = 'cl:10149'
$rows[2,1] = 'floating point'
$rows[2,2] = 'It is a numeric data type used to store numbers that may have a fractional component, like monetary values do (707.07, 0.7, 707.00).'
$rows[2,3] = 'float'
$rows[2,4] = 'https://dataled.academy/guides/data-types/'
$rows[2,5] = 'cl:10147'

$rows[3,0] = 'cl:10150'
$rows[3,1] = 'character'
$rows[3,2] = 'It is used to store a single letter, digit, punctuation mark, symbol, or blank space.'
$rows[3,3] = 'char'
$rows[3,4] = 'https://dataled.academy/guides/data-types/'
$rows[3,5] = 'cl:10147'

$rows[4,0] = 'cl:10151'
$rows[4,1] = 'string'
$rows[4,2] = 'It is a sequence of characters and the most commonly used data type to store text.'
$rows[4,3] = 'str'
$rows[4,4] = 'https://dataled.academy/guides/data-types/'
$rows[4,5] = 'cl:10147'

$rows[5,0] = 'cl:10152'
$rows[5,1] = 'boolean'
$rows[5,2] = 'It represents the values true and false.'
$rows[5,3] = 'bool'
$rows[5,4] = 'https://dataled.academy/guides/data-types/'
$rows[5,5] = 'cl:10147'

$rows[6,0] = 'cl:10153'
$rows[6,1] = 'enumeration type'
$rows[6,2] = 'It contains a small set of predefined unique values (also known as elements or enumerators) that can be compared and assigned to a variable of enumerated data type.'
$rows[6,3] = 'enum'
$rows[6,4] = 'https://dataled.academy/guides/data-types/'
$rows[6,5] = 'cl:10147'

$rows[7,0] = 'cl:10154'
$rows[7,1] = 'array'
$rows[7,2] = 'Also known as a list, an array is a data type that stores a number of elements in a specific order, typically all of the same type.'
$rows[7,3] = ''
$rows[7,4] = 'https://dataled.academy/guides/data-types/'
$rows[7,5] = 'cl:10147'

$rows[8,0] = 'cl:10155'
$rows[8,1] = 'date'
$rows[8,2] = 'It typically stores a date in the YYYY-MM-DD format (ISO 8601 syntax).'
$rows[8,3] = ''
$rows[8,4] = 'https://dataled.academy/guides/data-types/'
$rows[8,5] = 'cl:10147'

$rows[9,0] = 'cl:10156'
$rows[9,1] = 'time'
$rows[9,2] = 'Stores a time in the hh:mm:ss format.'
$rows[9,3] = ''
$rows[9,4] = 'https://dataled.academy/guides/data-types/'
$rows[9,5] = 'cl:10147'

$rows[10,0] = 'cl:10157'
$rows[10,1] = 'datetime'
$rows[10,2] = 'Stores a value containing both date and time together in the YYYY-MM-DD hh:mm:ss format.'
$rows[10,3] = ''
$rows[10,4] = 'https://dataled.academy/guides/data-types/'
$rows[10,5] = 'cl:10147'

$rows[11,0] = 'cl:10158'
$rows[11,1] = 'timestamp'
$rows[11,2] = 'Typically represented in Unix time, a timestamp represents the number of seconds that have elapsed since midnight (00:00:00 UTC), 1st January 1970.'
$rows[11,3] = ''
$rows[11,4] = 'https://dataled.academy/guides/data-types/'
$rows[11,5] = 'cl:10147'

$rows[12,0] = 'cl:10159'
$rows[12,1] = 'measurement scale'
$rows[12,2] = 'A measurement scale defines the level of measurement of a variable. It decides the statistical test type to be used. The mathematical nature of a variable or in other words, how a variable is measured is considered as the level of measurement.'
$rows[12,3] = ''
$rows[12,4] = ''
$rows[12,5] = ''

$rows[13,0] = 'cl:10160'
$rows[13,1] = 'nominal scale'
$rows[13,2] = 'Nominal scale is a naming scale, where variables are simply “named” or labeled, with no specific order.'
$rows[13,3] = ''
$rows[13,4] = ''
$rows[13,5] = 'cl:10160'

$rows[14,0] = 'cl:10161'
$rows[14,1] = 'ordinal scale'
$rows[14,2] = 'Ordinal scale has all its variables in a specific order, beyond just naming them. '
$rows[14,3] = ''
$rows[14,4] = ''
$rows[14,5] = 'cl:10160'

$rows[15,0] = 'cl:10162'
$rows[15,1] = 'interval scale'
$rows[15,2] = 'Interval scale offers labels, order, as well as, a specific interval between each of its variable options.'
$rows[15,3] = ''
$rows[15,4] = ''
$rows[15,5] = 'cl:10160'

$rows[16,0] = 'cl:10163'
$rows[16,1] = 'ratio scale'
$rows[16,2] = 'Ratio scale bears all the characteristics of an interval scale, in addition to that, it can also accommodate the value of “zero” on any of its variables.'
$rows[16,3] = ''
$rows[16,4] = ''
$rows[16,5] = 'cl:10160'

$startRow = 167
$endRow = 183
$rng = $ws.Range("A" + $startRow + ":F" + $endRow)
$rng.Value = $rows

Write-Output $ws.UsedRange.Address()
